$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.942.10'
$ws.Range('E2').Value = '  +0.18%  '

$ws.Range('D3').Value = '1.756.16'
$ws.Range('E3').Value = '  +0.12%  '

$ws.Range('E4').Value = '  +0.17%  '

$ws.Range('D5').Value = '''235.33'
$ws.Range('E5').Value = '  -1.78%  '

$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  +0.11%  '

$ws.Range('D7').Value = '''0.5212'
$ws.Range('E7').Value = '  +2.21%  '

$ws.Range('D8').Value = '''0.2737'
$ws.Range('E8').Value = '  -1.00%  '

$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = '''40.49'
$ws.Range('E9').Value = '  -4.19%  '

$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '''0.06162'
$ws.Range('E10').Value = '  -0.78%  '

$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').Value = '1.769.92'
$ws.Range('E11').Value = '  +0.93%  '

$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = '''0.07035'
$ws.Range('E12').Value = '  +0.81%  '

$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').Value = '''15.56'
$ws.Range('E13').Value = '  -1.61%  '

$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = '''0.6358'
$ws.Range('E14').Value = '  +3.53%  '

$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '''4.523'
$ws.Range('E15').Value = '  -0.32%  '

$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = '''77.69'
$ws.Range('E16').Value = '  +0.13%  '

$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('E17').Value = '  +0.25%  '

$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').Value = '''1.001'
$ws.Range('E18').Value = '  +0.13%  '

$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '25.957.63'
$ws.Range('E19').Value = '  +0.21%  '

$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = '''11.60'
$ws.Range('E20').Value = '  -0.70%  '

$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '''0.000006677'
$ws.Range('E21').Value = '  -3.68%  '

$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '1.995.94'
$ws.Range('E22').Value = '  +1.27%  '

$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = '''4.055'
$ws.Range('E23').Value = '  -0.83%  '

$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = '''8.474'
$ws.Range('E24').Value = '  +2.77%  '

$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D25').Value = '''5.174'
$ws.Range('E25').Value = '  -1.95%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '''138.98'
$ws.Range('E26').Value = '  +0.61%  '

$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').Value = '''1.502'
$ws.Range('E27').Value = '  +0.75%  '

$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = '''1.837'
$ws.Range('E28').Value = '  +0.53%  '

$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '''15.10'
$ws.Range('E29').Value = '  +0.05%  '

$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '''103.18'
$ws.Range('E30').Value = '  -0.55%  '

$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '''0.08376'
$ws.Range('E31').Value = '  +1.85%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '''3.662'
$ws.Range('E32').Value = '  -1.32%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '''3.422'
$ws.Range('E33').Value = '  -2.23%  '

$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '''0.04454'
$ws.Range('E34').Value = '  -2.05%  '

$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '''2.624'
$ws.Range('E35').Value = '  -0.77%  '

$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '''0.9939'
$ws.Range('E36').Value = '  -0.13%  '

$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '''0.6031'
$ws.Range('E37').Value = '  -1.77%  '

$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '''2.719'
$ws.Range('E38').Value = '  +0.00%  '

$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '''0.01586'
$ws.Range('E39').Value = '  +1.77%  '

$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '''1.955'
$ws.Range('E40').Value = '  +3.23%  '

$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = '''1.001'

$ws.Range('D42').Value = '''102.62'
$ws.Range('E42').Value = '  -1.44%  '

$ws.Range('D43').Value = '''0.3859'
$ws.Range('E43').Value = '  -0.65%  '

$ws.Range('D44').Value = '''0.7398'
$ws.Range('E44').Value = '  -0.42%  '

$ws.Range('D45').Value = '''4.904'
$ws.Range('E45').Value = '  -0.89%  '

$ws.Range('D46').Value = '''0.05512'
$ws.Range('E46').Value = '  +1.57%  '

$ws.Range('D47').Value = '''6.296'
$ws.Range('E47').Value = '  +4.41%  '

$ws.Range('D48').Value = '''0.1110'
$ws.Range('E48').Value = '  -0.67%  '

$ws.Range('D49').Value = '''30.07'
$ws.Range('E49').Value = '  -0.31%  '

$ws.Range('D50').Value = '''52.34'
$ws.Range('E50').Value = '  -1.08%  '

$ws.Range('D51').Value = '''1.002'

Write-Host "Applied cryptos update"
